$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Time_2 (F) and Time_3 (I) final times update, NAs removed from plot data
# Row 2
$ws.Range("F2").Value = 28.9
$ws.Range("I2").Value = 34.54

# Row 3
$ws.Range("F3").Value = 29.08
$ws.Range("I3").Value = 35.42

# Row 4
$ws.Range("F4").Value = 28.91
$ws.Range("I4").ClearContents()

# Row 5
$ws.Range("F5").Value = 29.2
$ws.Range("I5").Value = 35.31

# Row 6 (previously formula =22.33/1.0165)
$ws.Range("F6").ClearContents()
$ws.Range("I6").ClearContents()

# Row 7
$ws.Range("F7").Value = 28.911
$ws.Range("I7").Value = 34.8

# Row 8
$ws.Range("F8").Value = 28.89
$ws.Range("I8").Value = 34.67

# Row 9
$ws.Range("F9").ClearContents()
$ws.Range("I9").Value = 34.99

# Row 10
$ws.Range("F10").Value = 28.88
$ws.Range("I10").Value = 35.21

# Row 11
$ws.Range("F11").Value = 29.07
$ws.Range("I11").Value = 35.2

# Row 12
$ws.Range("F12").Value = 29.45
$ws.Range("I12").Value = 36.41

# Row 13
$ws.Range("F13").Value = 29.37
$ws.Range("I13").ClearContents()

# Row 14
$ws.Range("F14").Value = 28.98
$ws.Range("I14").Value = 34.92

# Row 15
$ws.Range("F15").Value = 29.41
$ws.Range("I15").Value = 36.58

# Update the selection to match the saved cursor position
$ws.Range("F23").Select()
